$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C30").Value = "어디에가요.mp3"
$ws.Range("C31").Value = "학교에가요.mp3"
$ws.Range("C32").Value = "도서관에서뭐해요.mp3"
$ws.Range("C33").Value = "도서관에서숙제를해요.mp3"
$ws.Range("C34").Value = "케빈팅팅씨어디에가요.mp3"
$ws.Range("C35").Value = "팅팅식당에가요식당에서친구하고밥을먹어요.mp3"
$ws.Range("C36").Value = "케빈아그래요.mp3"
$ws.Range("C37").Value = "팅팅케빈씨는오늘뭐해요.mp3"
$ws.Range("C38").Value = "케빈저는집에서쉬어요.mp3"
